$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update view: pane topLeftCell and active selection
$ws.Application.ActiveWindow.FreezePanes = $false
$ws.Range("E1").Select()
$ws.Application.ActiveWindow.FreezePanes = $true
$ws.Range("G25").Select()

# Update K column values (skill id numbers)
$ws.Range("K13").Value = 337389

$ws.Range("K23").Value = 326706
$ws.Range("K24").Value = 326704
$ws.Range("K25").Value = 326705
$ws.Range("K27").Value = 326710
$ws.Range("K28").Value = 327113
$ws.Range("K29").Value = 326693
$ws.Range("K30").Value = 326698
$ws.Range("K31").Value = 326699
$ws.Range("K32").Value = 326712

$ws.Range("K33").Value = 327119
$ws.Range("K35").Value = 327136
$ws.Range("K36").Value = 327101
$ws.Range("K38").Value = 327099
$ws.Range("K39").Value = 327118
$ws.Range("K40").Value = 327099
$ws.Range("K41").Value = 327126
$ws.Range("K42").Value = 327104

$ws.Range("K44").Value = 327083
$ws.Range("K45").Value = 327089
$ws.Range("K47").Value = 327088
$ws.Range("K48").Value = 327087
$ws.Range("K49").Value = 327170
$ws.Range("K50").Value = 327094
$ws.Range("K51").Value = 327171
$ws.Range("K52").Value = 327097
$ws.Range("K53").Value = 327100
